$d = $word.ActiveDocument

# Namespace-qualified package envelope used to splice in exact WordprocessingML
# for a given Range via Range.InsertXML (REPLACES that range's contents only).
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ------------------------------------------------------------------
# Change 1: "Nomor : 470/  68  /IX-2023/Ds." -> "Nomor : ${no_surat}"
#   Keep the "Nomor : " run untouched; collapse every run that carries
#   the old nomor-surat text into a single run holding the template
#   placeholder, reusing the first replaced run's rsid + rPr (minus lang).
# ------------------------------------------------------------------
$prefix1 = "Nomor : "
$para1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ($prefix1 + "*")) {
        $para1 = $p
        break
    }
}

if ($para1 -ne $null) {
    $start1 = $para1.Range.Start + $prefix1.Length
    $end1 = $para1.Range.End - 1
    $target1 = $d.Range($start1, $end1)

    $run1 = '<w:p><w:r w:rsidRPr="00903339"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="28"/></w:rPr><w:t>${no_surat}</w:t></w:r></w:p>'
    $target1.InsertXML($pkgOpen + $run1 + $pkgClose)
}

# ------------------------------------------------------------------
# Change 2: "01 September 2023" -> "${created_at}"
#   Keep the ": " run untouched; collapse the "01 September" / " 202" /
#   "3" runs into a single run holding the template placeholder, reusing
#   the first replaced run's (lang-less) rPr.
# ------------------------------------------------------------------
$prefix2 = ": "
$para2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ($prefix2 + "01 September*")) {
        $para2 = $p
        break
    }
}

if ($para2 -ne $null) {
    $start2 = $para2.Range.Start + $prefix2.Length
    $end2 = $para2.Range.End - 1
    $target2 = $d.Range($start2, $end2)

    $run2 = '<w:p><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>${created_at}</w:t></w:r></w:p>'
    $target2.InsertXML($pkgOpen + $run2 + $pkgClose)
}

Write-Output "Para1 after: $($d.Paragraphs | Where-Object { $_.Range.Text -like 'Nomor *' } | Select-Object -First 1 -ExpandProperty Range | Select-Object -ExpandProperty Text)"
